$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column E ("Publon Profile" onward shifts right by one),
# shifting all data in columns E:I to F:J.
$ws.Columns("E:E").Insert()

# New header for the inserted column.
$ws.Range("E1").Value = "Research Gate"

# Match the author's selection (E1) left after the edit.
$ws.Range("E1").Select()

# Column widths (D..K) as set by the author after the insert.
$ws.Columns("D:D").ColumnWidth = 56.333333333333336
$ws.Columns("E:E").ColumnWidth = 43.333333333333336
$ws.Columns("F:F").ColumnWidth = 30.833333333333332
$ws.Columns("G:G").ColumnWidth = 19.666666666666668
$ws.Columns("H:H").ColumnWidth = 21.166666666666668
$ws.Columns("I:I").ColumnWidth = 24.333333333333332
$ws.Columns("J:K").ColumnWidth = 18.166666666666668
